$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42: E42 text change from {'any', 'str'} to {'str', 'any'}
$ws.Range("E42").Value = "{'str', 'any'}"

# Row 43: E43 any -> str ; F43 Loss -> Neutral (style changes from red fill to orange fill)
$ws.Range("E43").Value = "str"
$ws.Range("F43").Value = "Neutral"
$ws.Range("F43").Interior.Color = $ws.Range("F42").Interior.Color

# Row 74: D74 7 -> 6 (PyType Wins count)
$ws.Range("D74").Value = 6

# Insert a new row 76, copying formatting from row 75 (the current last row)
# so every cell A76:F76 carries the same style (s="2") as the rest of the table.
$ws.Range("A75:F75").Copy()
$ws.Range("A76:F76").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 75: repurpose as "Scalpel Accuracy:" row, clear old "Accuracy over PyType" content
$ws.Range("C75").Value = "Scalpel Accuracy:"
$ws.Range("D75").Value = 1100
$ws.Range("E75").Value = ""
$ws.Range("F75").Value = ""

# Row 76: new row with "Accuracy over PyType" moved here with updated value
$ws.Range("A76").Value = ""
$ws.Range("B76").Value = ""
$ws.Range("C76").Value = ""
$ws.Range("D76").Value = ""
$ws.Range("E76").Value = "Accuracy over PyType"
$ws.Range("F76").Value = 50
